$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure all touched cells keep their original text formatting (many values
# look numeric, e.g. "22.50" or "0.06534", but must stay literal text).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.563.37"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.62%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.884.11"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.26%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "245.74"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.87%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4740"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2892"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06534"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.02%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "22.50"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.76%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.7743"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +4.88%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "100.82"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +4.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07819"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +0.04%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.885.18"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "5.250"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "285.63"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.556.95"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.61%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007527"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -0.25%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.131.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -0.32%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.351"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.01%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "6.444"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.94%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "9.181"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -0.23%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "163.44"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.77%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.12"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +0.65%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.912"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.09%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.341"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.09705"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.33%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.500"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.50%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.253"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.16%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.190"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +0.09%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04846"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.32%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.129"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +0.00%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6980"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.757"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  +1.19%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.01907"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +0.82%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +3.44%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "76.00"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.41%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "6.295"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.46%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.984"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -0.54%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.4258"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.62%  "
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8348"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.11%  "
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = "EnergySwap"
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "9.940"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.99%  "
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = "Quant"
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.60"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "7.033"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.13%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "35.19"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -1.26%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05774"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  +0.24%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "Maker"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "891.62"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.24%  "
